$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = @(
  @('Examtaker34245','Automation34245','examtakerautomation34245@gmail.com','34245'),
  @('Examtaker67841','Automation67841','examtakerautomation67841@gmail.com','67841'),
  @('Examtaker78837','Automation78837','examtakerautomation78837@gmail.com','78837'),
  @('Examtaker18080','Automation18080','examtakerautomation18080@gmail.com','18080'),
  @('Examtaker58765','Automation58765','examtakerautomation58765@gmail.com','58765'),
  @('Examtaker33514','Automation33514','examtakerautomation33514@gmail.com','33514'),
  @('Examtaker38291','Automation38291','examtakerautomation38291@gmail.com','38291'),
  @('Examtaker60430','Automation60430','examtakerautomation60430@gmail.com','60430'),
  @('Examtaker02610','Automation02610','examtakerautomation02610@gmail.com','02610'),
  @('Examtaker95821','Automation95821','examtakerautomation95821@gmail.com','95821'),
  @('Examtaker81159','Automation81159','examtakerautomation81159@gmail.com','81159'),
  @('Examtaker98224','Automation98224','examtakerautomation98224@gmail.com','98224'),
  @('Examtaker56331','Automation56331','examtakerautomation56331@gmail.com','56331'),
  @('Examtaker91917','Automation91917','examtakerautomation91917@gmail.com','91917'),
  @('Examtaker35284','Automation35284','examtakerautomation35284@gmail.com','35284'),
  @('Examtaker52060','Automation52060','examtakerautomation52060@gmail.com','52060'),
  @('Examtaker71912','Automation71912','examtakerautomation71912@gmail.com','71912'),
  @('Examtaker75613','Automation75613','examtakerautomation75613@gmail.com','75613'),
  @('Examtaker63210','Automation63210','examtakerautomation63210@gmail.com','63210'),
  @('Examtaker83857','Automation83857','examtakerautomation83857@gmail.com','83857'),
  @('Examtaker54033','Automation54033','examtakerautomation54033@gmail.com','54033'),
  @('Examtaker97377','Automation97377','examtakerautomation97377@gmail.com','97377'),
  @('Examtaker08744','Automation08744','examtakerautomation08744@gmail.com','08744'),
  @('Examtaker84660','Automation84660','examtakerautomation84660@gmail.com','84660'),
  @('Examtaker30569','Automation30569','examtakerautomation30569@gmail.com','30569'),
  @('Examtaker00354','Automation00354','examtakerautomation00354@gmail.com','00354'),
  @('Examtaker95822','Automation95822','examtakerautomation95822@gmail.com','95822'),
  @('Examtaker36463','Automation36463','examtakerautomation36463@gmail.com','36463'),
  @('Examtaker57591','Automation57591','examtakerautomation57591@gmail.com','57591'),
  @('Examtaker27719','Automation27719','examtakerautomation27719@gmail.com','27719')
)
$r = 2
foreach ($rec in $data) {
    $ws.Cells.Item($r,1).Value = $rec[0]
    $ws.Cells.Item($r,2).Value = $rec[1]
    $ws.Cells.Item($r,3).Value = $rec[2]
    $dcell = $ws.Cells.Item($r,4)
    $dcell.NumberFormat = "@"
    $dcell.Value = $rec[3]
    $dcell.Style = "Normal"
    $r = $r + 1
}
